$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55
$ws.Range("A55").Value = 111898336
$ws.Range("B55").Value = 89405
$ws.Range("C55").Value = "Ovaliderad"
$ws.Range("D55").Value = "NT"
$ws.Range("E55").Value = 1202
$ws.Range("F55").Value = "Ullticka"
$ws.Range("G55").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H55").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P55").Value = "Lill-Ångeskogen 3 km SV om Lena kyrka, Upl"
$ws.Range("Q55").Value = 650105.085176448
$ws.Range("R55").Value = 6654011.298884101
$ws.Range("S55").Value = 10
$ws.Range("T55").Value = "Uppsala"
$ws.Range("U55").Value = "Uppsala"
$ws.Range("V55").Value = "Uppland"
$ws.Range("W55").Value = "Lena"
$ws.Range("Y55").NumberFormat = "@"
$ws.Range("Y55").Value = "2023-09-03"
$ws.Range("Z55").NumberFormat = "@"
$ws.Range("Z55").Value = "00:00"
$ws.Range("AA55").NumberFormat = "@"
$ws.Range("AA55").Value = "2023-09-03"
$ws.Range("AB55").NumberFormat = "@"
$ws.Range("AB55").Value = "00:00"
$ws.Range("AD55").Value = $false
$ws.Range("AE55").Value = $false
$ws.Range("AG55").Value = $false
$ws.Range("AH55").Value = "Ängsblandskog"
$ws.Range("AJ55").Value = "gran"
$ws.Range("AK55").Value = "Picea abies"
$ws.Range("AM55").Value = "Liggande död trädstam, utan markontakt"
$ws.Range("AO55").Value = "Horizontal, dead without ground contact # Picea abies"
$ws.Range("AW55").Value = "Thorleif Joelson"
$ws.Range("AX55").Value = "Thorleif Joelson, Henry Åkerström"
$ws.Range("AY55").Value = "Naturskyddsföreningen Uppsala, skogsgruppexkursion"

# Row 56
$ws.Range("A56").Value = 111898889
$ws.Range("B56").Value = 98535
$ws.Range("C56").Value = "Ovaliderad"
$ws.Range("D56").Value = "LC"
$ws.Range("E56").Value = 222498
$ws.Range("F56").Value = "Blåsippa"
$ws.Range("G56").Value = "Hepatica nobilis"
$ws.Range("H56").Value = "Schreb."
$ws.Range("K56").Value = "fullt utvecklade blad"
$ws.Range("P56").Value = "Lill-Ångeskogen 3 km SV om Lena kyrka, Upl"
$ws.Range("Q56").Value = 650135.0421630922
$ws.Range("R56").Value = 6654002.501842719
$ws.Range("S56").Value = 10
$ws.Range("T56").Value = "Uppsala"
$ws.Range("U56").Value = "Uppsala"
$ws.Range("V56").Value = "Uppland"
$ws.Range("W56").Value = "Lena"
$ws.Range("Y56").NumberFormat = "@"
$ws.Range("Y56").Value = "2023-09-03"
$ws.Range("Z56").NumberFormat = "@"
$ws.Range("Z56").Value = "00:00"
$ws.Range("AA56").NumberFormat = "@"
$ws.Range("AA56").Value = "2023-09-03"
$ws.Range("AB56").NumberFormat = "@"
$ws.Range("AB56").Value = "00:00"
$ws.Range("AD56").Value = $false
$ws.Range("AE56").Value = $false
$ws.Range("AG56").Value = $false
$ws.Range("AH56").Value = "Ängsbarrskog"
$ws.Range("AI56").Value = "Ungskog"
$ws.Range("AW56").Value = "Thorleif Joelson"
$ws.Range("AX56").Value = "Thorleif Joelson, Henry Åkerström"
$ws.Range("AY56").Value = "Naturskyddsföreningen Uppsala, skogsgruppexkursion"

# Row 57
$ws.Range("A57").Value = 111898507
$ws.Range("B57").Value = 89845
$ws.Range("C57").Value = "Ovaliderad"
$ws.Range("D57").Value = "VU"
$ws.Range("E57").Value = 1209
$ws.Range("F57").Value = "Rynkskinn"
$ws.Range("G57").Value = "Phlebia centrifuga"
$ws.Range("H57").Value = "P.Karst."
$ws.Range("P57").Value = "Lill-Ångeskogen 3 km SV om Lena kyrka, Upl"
$ws.Range("Q57").Value = 650086.8716060545
$ws.Range("R57").Value = 6654015.064976334
$ws.Range("S57").Value = 10
$ws.Range("T57").Value = "Uppsala"
$ws.Range("U57").Value = "Uppsala"
$ws.Range("V57").Value = "Uppland"
$ws.Range("W57").Value = "Lena"
$ws.Range("Y57").NumberFormat = "@"
$ws.Range("Y57").Value = "2023-09-03"
$ws.Range("Z57").NumberFormat = "@"
$ws.Range("Z57").Value = "00:00"
$ws.Range("AA57").NumberFormat = "@"
$ws.Range("AA57").Value = "2023-09-03"
$ws.Range("AB57").NumberFormat = "@"
$ws.Range("AB57").Value = "00:00"
$ws.Range("AD57").Value = $false
$ws.Range("AE57").Value = $false
$ws.Range("AG57").Value = $false
$ws.Range("AH57").Value = "Ängsblandskog"
$ws.Range("AJ57").Value = "gran"
$ws.Range("AK57").Value = "Picea abies"
$ws.Range("AM57").Value = "Liggande död trädstam, utan markontakt"
$ws.Range("AO57").Value = "Horizontal, dead without ground contact # Picea abies"
$ws.Range("AW57").Value = "Thorleif Joelson"
$ws.Range("AX57").Value = "Thorleif Joelson, Henry Åkerström"
$ws.Range("AY57").Value = "Naturskyddsföreningen Uppsala, skogsgruppexkursion"

# Row 58
$ws.Range("A58").Value = 111898660
$ws.Range("B58").Value = 100532
$ws.Range("C58").Value = "Ovaliderad"
$ws.Range("D58").Value = "CR"
$ws.Range("E58").Value = 223246
$ws.Range("F58").Value = "Skogsalm"
$ws.Range("G58").Value = "Ulmus glabra"
$ws.Range("H58").Value = "Huds."
$ws.Range("P58").Value = "Lill-Ångeskogen 3 km SV om Lena kyrka, Upl"
$ws.Range("Q58").Value = 650054.1336129439
$ws.Range("R58").Value = 6654018.240072312
$ws.Range("S58").Value = 10
$ws.Range("T58").Value = "Uppsala"
$ws.Range("U58").Value = "Uppsala"
$ws.Range("V58").Value = "Uppland"
$ws.Range("W58").Value = "Lena"
$ws.Range("Y58").NumberFormat = "@"
$ws.Range("Y58").Value = "2023-09-03"
$ws.Range("Z58").NumberFormat = "@"
$ws.Range("Z58").Value = "00:00"
$ws.Range("AA58").NumberFormat = "@"
$ws.Range("AA58").Value = "2023-09-03"
$ws.Range("AB58").NumberFormat = "@"
$ws.Range("AB58").Value = "00:00"
$ws.Range("AC58").Value = "Stammens omkrets i brösthöjd: 64 cm"
$ws.Range("AD58").Value = $false
$ws.Range("AE58").Value = $false
$ws.Range("AG58").Value = $false
$ws.Range("AH58").Value = "Ängsblandskog"
$ws.Range("AW58").Value = "Thorleif Joelson"
$ws.Range("AX58").Value = "Thorleif Joelson, Henry Åkerström"
$ws.Range("AY58").Value = "Naturskyddsföreningen Uppsala, skogsgruppexkursion"

# Row 59
$ws.Range("A59").Value = 111898191
$ws.Range("B59").Value = 90332
$ws.Range("C59").Value = "Ovaliderad"
$ws.Range("D59").Value = "LC"
$ws.Range("E59").Value = 4769
$ws.Range("F59").Value = "Svavelriska"
$ws.Range("G59").Value = "Lactarius scrobiculatus"
$ws.Range("H59").Value = "(Scop.:Fr.) Fr."
$ws.Range("I59").NumberFormat = "@"
$ws.Range("I59").Value = "2"
$ws.Range("J59").Value = "fruktkroppar"
$ws.Range("P59").Value = "Lill-Ångeskogen 3 km SV om Lena kyrka, Upl"
$ws.Range("Q59").Value = 650135.0421630922
$ws.Range("R59").Value = 6654002.501842719
$ws.Range("S59").Value = 10
$ws.Range("T59").Value = "Uppsala"
$ws.Range("U59").Value = "Uppsala"
$ws.Range("V59").Value = "Uppland"
$ws.Range("W59").Value = "Lena"
$ws.Range("Y59").NumberFormat = "@"
$ws.Range("Y59").Value = "2023-09-03"
$ws.Range("Z59").NumberFormat = "@"
$ws.Range("Z59").Value = "00:00"
$ws.Range("AA59").NumberFormat = "@"
$ws.Range("AA59").Value = "2023-09-03"
$ws.Range("AB59").NumberFormat = "@"
$ws.Range("AB59").Value = "00:00"
$ws.Range("AD59").Value = $false
$ws.Range("AE59").Value = $false
$ws.Range("AG59").Value = $false
$ws.Range("AH59").Value = "Ängsbarrskog"
$ws.Range("AI59").Value = "Ungskog"
$ws.Range("AW59").Value = "Thorleif Joelson"
$ws.Range("AX59").Value = "Thorleif Joelson, Henry Åkerström"
$ws.Range("AY59").Value = "Naturskyddsföreningen Uppsala, skogsgruppexkursion"

# Row 60
$ws.Range("A60").Value = 111911660
$ws.Range("B60").Value = 96348
$ws.Range("C60").Value = "Ovaliderad"
$ws.Range("D60").Value = "VU"
$ws.Range("E60").Value = 220787
$ws.Range("F60").Value = "Knärot"
$ws.Range("G60").Value = "Goodyera repens"
$ws.Range("H60").Value = "(L.) R. Br."
$ws.Range("I60").NumberFormat = "@"
$ws.Range("I60").Value = "19"
$ws.Range("J60").Value = "plantor/tuvor"
$ws.Range("P60").Value = "Lill-Ångeskogen, 3 km SV om Lena kyrka, Upl"
$ws.Range("Q60").Value = 650026.652882754
$ws.Range("R60").Value = 6654299.07778531
$ws.Range("S60").Value = 10
$ws.Range("T60").Value = "Uppsala"
$ws.Range("U60").Value = "Uppsala"
$ws.Range("V60").Value = "Uppland"
$ws.Range("W60").Value = "Lena"
$ws.Range("Y60").NumberFormat = "@"
$ws.Range("Y60").Value = "2023-09-03"
$ws.Range("Z60").NumberFormat = "@"
$ws.Range("Z60").Value = "00:00"
$ws.Range("AA60").NumberFormat = "@"
$ws.Range("AA60").Value = "2023-09-03"
$ws.Range("AB60").NumberFormat = "@"
$ws.Range("AB60").Value = "00:00"
$ws.Range("AD60").Value = $false
$ws.Range("AE60").Value = $false
$ws.Range("AG60").Value = $false
$ws.Range("AH60").Value = "Barrskog på blockmark"
$ws.Range("AW60").Value = "Thorleif Joelson"
$ws.Range("AX60").Value = "Thorleif Joelson, Henry Åkerström"
$ws.Range("AY60").Value = "Naturskyddsföreningen Uppsala, skogsgruppexkursion"

# Row 61
$ws.Range("A61").Value = 111911698
$ws.Range("B61").Value = 96348
$ws.Range("C61").Value = "Ovaliderad"
$ws.Range("D61").Value = "VU"
$ws.Range("E61").Value = 220787
$ws.Range("F61").Value = "Knärot"
$ws.Range("G61").Value = "Goodyera repens"
$ws.Range("H61").Value = "(L.) R. Br."
$ws.Range("I61").NumberFormat = "@"
$ws.Range("I61").Value = "16"
$ws.Range("J61").Value = "plantor/tuvor"
$ws.Range("P61").Value = "Lill-Ångeskogen, 3 km SV om Lena kyrka, Upl"
$ws.Range("Q61").Value = 650032.9755174413
$ws.Range("R61").Value = 6654279.303373625
$ws.Range("S61").Value = 10
$ws.Range("T61").Value = "Uppsala"
$ws.Range("U61").Value = "Uppsala"
$ws.Range("V61").Value = "Uppland"
$ws.Range("W61").Value = "Lena"
$ws.Range("Y61").NumberFormat = "@"
$ws.Range("Y61").Value = "2023-09-03"
$ws.Range("Z61").NumberFormat = "@"
$ws.Range("Z61").Value = "00:00"
$ws.Range("AA61").NumberFormat = "@"
$ws.Range("AA61").Value = "2023-09-03"
$ws.Range("AB61").NumberFormat = "@"
$ws.Range("AB61").Value = "00:00"
$ws.Range("AD61").Value = $false
$ws.Range("AE61").Value = $false
$ws.Range("AG61").Value = $false
$ws.Range("AH61").Value = "Barrskog på blockmark"
$ws.Range("AW61").Value = "Thorleif Joelson"
$ws.Range("AX61").Value = "Thorleif Joelson, Henry Åkerström"
$ws.Range("AY61").Value = "Naturskyddsföreningen Uppsala, skogsgruppexkursion"

Write-Host "Added rows 55-61"